$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.886.68"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.806.66"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.76"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +3.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3715"
$ws.Range("E8").Value = "  -1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07381"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8736"
$ws.Range("E10").Value = "  -1.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.44"
$ws.Range("E11").Value = "  -2.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.830.63"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.379"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.91"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.493"
$ws.Range("E15").Value = "  -3.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07031"
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008728"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.70"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.896.66"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.306"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.034.69"
$ws.Range("E24").Value = "  -0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("E25").Value = "  -3.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.46"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.156"
$ws.Range("E28").Value = "  -6.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.287"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.98"
$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08933"
$ws.Range("E31").Value = "  +0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7602"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.157"
$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.457"
$ws.Range("E34").Value = "  -3.34%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05258"
$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("E40").Value = "  +2.11%  "

$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.379"
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5289"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("E44").Value = "  -2.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.528"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4994"
$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.31"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.99"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.667"
$ws.Range("E50").Value = "  -1.60%  "
